# Updates the crypto price/volume table to the latest scraped values.
# (matches the data refresh performed by the GitHub Actions job)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.688.73'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.633.87'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.19'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.498'
$ws.Range('E6').Value = '  -1.54%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.95'
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').Value = '1.860.40'
$ws.Range('E12').Value = '  -0.92%  '
$ws.Range('D13').Value = '1.636.35'
$ws.Range('E13').Value = '  -1.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.11'
$ws.Range('E14').Value = '  -2.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.522'
$ws.Range('E15').Value = '  -2.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.98'
$ws.Range('E16').Value = '  -2.32%  '
$ws.Range('D17').Value = '26.672.55'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').Value = '0.0₃0722'
$ws.Range('E18').Value = '  -3.06%  '
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '211.01'
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('E21').Value = '  -1.66%  '
$ws.Range('E22').Value = '  -2.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.31'
$ws.Range('E23').Value = '  -8.93%  '
$ws.Range('E24').Value = '  -3.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.75'
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  -2.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.00'
$ws.Range('E28').Value = '  -2.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.54'
$ws.Range('E29').Value = '  -1.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0501'
$ws.Range('E30').Value = '  -3.88%  '
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.36'
$ws.Range('E32').Value = '  +0.27%  '
$ws.Range('E33').Value = '  -2.81%  '
$ws.Range('D34').Value = '1.260.83'
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('E35').Value = '  -2.36%  '
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('E37').Value = '  -3.50%  '
$ws.Range('E38').Value = '  -3.35%  '
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.802'
$ws.Range('E40').Value = '  -3.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.799'
$ws.Range('E41').Value = '  -2.42%  '
$ws.Range('E42').Value = '  -4.28%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.25'
$ws.Range('E43').Value = '  -3.08%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.769.80'
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.30'
$ws.Range('E45').Value = '  -0.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.76'
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('E47').Value = '  -3.49%  '
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0955'
$ws.Range('E51').Value = '  -2.58%  '
